# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
# Corrects mismatched fixture rows (each row's id in col B plus all the
# match data in columns F:AC had been written against the wrong fixture).
# The fix re-pairs the per-row payload (F:AC, and the id in B) with the
# correct row, leaving A (running index), C/D (Div), and E (Date) as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ecuador LigaPro Serie A")

# Columns used: B=id, F=HomeTeam, G=AwayTeam, H=FTHG, I=FTAG, J=FTR,
# K=oddH_op, L=oddD_op, M=oddA_op, N=oddH, O=oddD, P=oddA, Q=Ah,
# R=oddAHH, S=oddAHA, T=AhOU, U=oddAHOver, V=oddAHUnder, W=PLH, X=PLD,
# Y=PLA, Z=PL_Ahh, AA=PL_Aha, AB=PL_AhOver, AC=PL_AhUnder
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowData($row) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value2 = $data[$c]
    }
}

# Snapshot current (pre-edit) payloads for every affected row first, since
# several rows are written based on each other's original values.
$r130 = Get-RowData 130
$r132 = Get-RowData 132
$r133 = Get-RowData 133

$r142 = Get-RowData 142
$r143 = Get-RowData 143
$r144 = Get-RowData 144
$r145 = Get-RowData 145

$r214 = Get-RowData 214
$r215 = Get-RowData 215

# Block 1 (rows 130, 132, 133): 3-way rotation.
Set-RowData 130 $r132
Set-RowData 132 $r133
Set-RowData 133 $r130

# Block 2 (rows 142/143 and 144/145): pairwise swaps.
Set-RowData 142 $r143
Set-RowData 143 $r142
Set-RowData 144 $r145
Set-RowData 145 $r144

# Block 3 (rows 214, 215): swap.
Set-RowData 214 $r215
Set-RowData 215 $r214
